$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '62.818.20'
$ws.Cells.Item(2, 5).Value = '  -1.94%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.673.99'
$ws.Cells.Item(3, 5).Value = '  -2.39%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '551.63'
$ws.Cells.Item(5, 5).Value = '  -3.19%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '157.82'
$ws.Cells.Item(6, 5).Value = '  -0.69%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.06%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.592'
$ws.Cells.Item(8, 5).Value = '  -1.03%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.105'
$ws.Cells.Item(9, 5).Value = '  -3.57%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -1.67%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -4.23%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '5.28'
$ws.Cells.Item(12, 5).Value = '  -7.67%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.145.24'
$ws.Cells.Item(13, 5).Value = '  -2.34%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -2.17%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '62.690.54'
$ws.Cells.Item(15, 5).Value = '  -1.48%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -2.73%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.671.81'
$ws.Cells.Item(17, 5).Value = '  -2.61%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '11.85'
$ws.Cells.Item(18, 5).Value = '  -1.96%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '4.60'
$ws.Cells.Item(19, 5).Value = '  -4.26%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '344.52'
$ws.Cells.Item(20, 5).Value = '  -2.78%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '6.27'
$ws.Cells.Item(21, 5).Value = '  -4.95%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.04%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -3.49%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '63.18'
$ws.Cells.Item(24, 5).Value = '  -1.82%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.82%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.33%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -3.14%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '0.0₃0858'
$ws.Cells.Item(28, 5).Value = '  -5.85%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '1.37'
$ws.Cells.Item(29, 5).Value = '  +1.50%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '7.24'
$ws.Cells.Item(30, 5).Value = '  +0.18%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -1.59%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '165.99'
$ws.Cells.Item(32, 5).Value = '  +1.05%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.23%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  +0.01%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(35, 4).Value = '4.85'
$ws.Cells.Item(35, 5).Value = '  -1.10%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '19.47'
$ws.Cells.Item(36, 5).Value = '  -2.72%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.78'
$ws.Cells.Item(37, 5).Value = '  -1.69%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '346.45'
$ws.Cells.Item(38, 5).Value = '  -1.09%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.955'
$ws.Cells.Item(39, 5).Value = '  -3.70%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '6.26'
$ws.Cells.Item(40, 5).Value = '  -0.90%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -3.29%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '38.24'
$ws.Cells.Item(42, 5).Value = '  -0.86%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(43, 4).Value = '20.78'
$ws.Cells.Item(43, 5).Value = '  -5.79%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = '20.30'
$ws.Cells.Item(44, 5).Value = '  -4.05%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.0563'
$ws.Cells.Item(45, 5).Value = '  -3.76%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -1.69%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.998'
$ws.Cells.Item(47, 5).Value = '  +0.01%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.0971'
$ws.Cells.Item(49, 5).Value = '  -3.32%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50, 4).Value = '0.0241'
$ws.Cells.Item(50, 5).Value = '  -2.85%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '128.62'
$ws.Cells.Item(51, 5).Value = '  -4.52%  '
